$d = $word.ActiveDocument

# --- Step 1: strip the _GoBack bookmark from the end of the "R: FEITO"
# paragraph (it moves further down, into the reworded answer paragraph).
# We rewrite that paragraph in place, preserving all of its original
# paragraph/run attributes, just dropping the bookmark.
$pFeito = $d.Paragraphs.Item(10)
$rFeito = $pFeito.Range
$rFeito.Collapse(1)
$xmlFeito = '<w:p xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="5385A14A" w14:textId="724DCC0E" w:rsidR="002153C2" w:rsidRPr="002153C2" w:rsidRDefault="002153C2" w:rsidP="006417F8"><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="002153C2"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>R: FEITO</w:t></w:r></w:p>'
$rFeito.InsertXML($xmlFeito)

# --- Step 2: replace the old "R:  Sera seguido..." paragraph (currently
# paragraph 11, right after "R: FEITO") with the two new paragraphs in the
# new order: the "2a (peso 2,0)..." question first, then the reworded,
# fully-bold answer (with the underlined model name and the relocated
# _GoBack bookmark).
$pOldAnswer = $d.Paragraphs.Item(11)
$rOldAnswer = $pOldAnswer.Range
$rOldAnswer.Collapse(1)
$xmlNew = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>2</w:t></w:r><w:r><w:t xml:space="preserve">ª (peso 2,0): </w:t></w:r><w:r><w:t>Considerando as condições do projeto, qual o Ciclo de Vida de produção de software que você seguirá? Justifique a sua resposta</w:t></w:r><w:r><w:t xml:space="preserve"> em um documento WORD chamado </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>PropostaSIColetaLixo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, o qual você vai subir no GITHUB, na pasta </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Documentacao</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">R:  Será seguido o </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Modelo de Prototipação evolucionária</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>, pois este modelo exige flexibilidade quanto à possibilidade de mud</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">ar de escopo, além do dono da empresa exigir uma entrega o quanto antes algum recurso, onde neste modelo podemos fazemos um </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>mockup</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> e já sair testando</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">. </w:t></w:r></w:p>'
$rOldAnswer.InsertXML($xmlNew)

# --- Step 3: the original "2a (peso 2,0)..." question paragraph is now a
# duplicate sitting right after our freshly-inserted replacement content
# (paragraph 13); along with what is left of the old answer paragraph
# (already overwritten above) nothing else needs removing except that
# duplicate question paragraph.
$dup = $d.Paragraphs.Item(13)
$dup.Range.Delete()
